# Apply crypto price/volume updates per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.624.55"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "2.285.33"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'96.55"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").Value = "'266.82"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "'0.610"
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("D10").Value = "'45.94"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").Value = "'0.0931"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "'7.81"
$ws.Range("E12").Value = "  -2.94%  "
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").Value = "2.629.48"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "'15.13"
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("D16").Value = "'0.835"
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("D17").Value = "2.287.45"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "43.613.59"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").Value = "'0.0000108"
$ws.Range("E19").Value = "  +2.32%  "
$ws.Range("D20").Value = "'6.21"
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("D21").Value = "'72.08"
$ws.Range("E21").Value = "  +1.34%  "
$ws.Range("D22").Value = "'2.45"
$ws.Range("E22").Value = "  +6.59%  "
$ws.Range("D23").Value = "'232.82"
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("D24").Value = "'9.20"
$ws.Range("E24").Value = "  -7.45%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").Value = "'2.50"
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("E28").Value = "  +2.26%  "
$ws.Range("D29").Value = "'40.46"
$ws.Range("E29").Value = "  +2.66%  "
$ws.Range("D30").Value = "'2.22"
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("D31").Value = "'176.68"
$ws.Range("E31").Value = "  +2.41%  "
$ws.Range("D32").Value = "'21.83"
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("D33").Value = "'0.0891"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("D34").Value = "'5.36"
$ws.Range("E34").Value = "  -3.71%  "
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("D36").Value = "'0.109"
$ws.Range("E36").Value = "  -1.95%  "
$ws.Range("D37").Value = "'0.0355"
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("D38").Value = "'4.32"
$ws.Range("E38").Value = "  -3.72%  "
$ws.Range("D39").Value = "'3.39"
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("D41").Value = "'2.31"
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("D43").Value = "'1.35"
$ws.Range("E43").Value = "  +1.47%  "
$ws.Range("D44").Value = "'64.92"
$ws.Range("E44").Value = "  +5.85%  "
$ws.Range("D45").Value = "'8.78"
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").Value = "'5.20"
$ws.Range("E46").Value = "  -4.99%  "
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'97.05"
$ws.Range("E48").Value = "  -3.14%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").Value = "'1.19"
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("D50").Value = "'0.434"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").Value = "2.509.10"
$ws.Range("E51").Value = "  -0.39%  "
